$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Amazon-1")

# Assumed growth-rate inputs (column I) -- percentage formatted
$ws.Range("I3").Value = 0.31
$ws.Range("I3").NumberFormat = "0%"
$ws.Range("I6").Value = 0.22
$ws.Range("I6").NumberFormat = "0%"
$ws.Range("I8").Value = 0.35
$ws.Range("I8").NumberFormat = "0%"

# Row 3 projected values now reference the assumed growth rate in I3
$ws.Range("L3").Formula = "=H3 *(1 + I3)"
$ws.Range("L3").ClearFormats()
$ws.Range("M3").Formula = "=L3 *(1+ I3)"
$ws.Range("M3").ClearFormats()

# Row 6 projected values now reference the assumed growth rate in I6
$ws.Range("L6").Formula = "=H6*(1 + I6)"
$ws.Range("L6").ClearFormats()
$ws.Range("M6").Formula = "=L6*(1 + I6)"
$ws.Range("M6").ClearFormats()

# New row 9 projected values based on the assumed growth rate in I8
$ws.Range("L9").Formula = "=H8*(1+I8)"
$ws.Range("L9").NumberFormat = "#,##0"
$ws.Range("M9").Formula = "=L9*(1 + I8)"
$ws.Range("M9").ClearFormats()

# Move the active selection to L11
$ws.Range("L11").Select()

$wb.Save()
